$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "28.362.71"
$cell.ClearFormats()
$ws.Cells.Item(2, 5).Value = "  +4.15%  "

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.595.05"
$cell.ClearFormats()
$ws.Cells.Item(3, 5).Value = "  +1.96%  "

$ws.Cells.Item(4, 5).Value = "  +0.01%  "

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "214.40"
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = "  +1.78%  "

$ws.Cells.Item(6, 5).Value = "  +1.33%  "

$ws.Cells.Item(7, 5).Value = "  +0.00%  "

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "24.01"
$cell.ClearFormats()
$ws.Cells.Item(8, 5).Value = "  +8.47%  "

$ws.Cells.Item(9, 5).Value = "  +0.68%  "

$ws.Cells.Item(10, 5).Value = "  +0.94%  "

$ws.Cells.Item(11, 5).Value = "  +2.00%  "

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.823.76"
$cell.ClearFormats()
$ws.Cells.Item(12, 5).Value = "  +2.02%  "

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.597.86"
$cell.ClearFormats()
$ws.Cells.Item(13, 5).Value = "  +1.93%  "

$ws.Cells.Item(14, 5).Value = "  +2.76%  "

$ws.Cells.Item(15, 5).Value = "  -0.10%  "

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "28.378.47"
$cell.ClearFormats()
$ws.Cells.Item(16, 5).Value = "  +4.38%  "

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "63.15"
$cell.ClearFormats()
$ws.Cells.Item(17, 5).Value = "  +2.03%  "

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "227.86"
$cell.ClearFormats()
$ws.Cells.Item(18, 5).Value = "  +4.66%  "

$ws.Cells.Item(19, 5).Value = "  +1.37%  "

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.49"
$cell.ClearFormats()
$ws.Cells.Item(20, 5).Value = "  +0.52%  "

$ws.Cells.Item(21, 5).Value = "  -0.10%  "

$ws.Cells.Item(22, 5).Value = "  -0.81%  "

$ws.Cells.Item(23, 5).Value = "  -0.42%  "

$ws.Cells.Item(24, 5).Value = "  +0.62%  "

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "151.58"
$cell.ClearFormats()
$ws.Cells.Item(25, 5).Value = "  -0.03%  "

$ws.Cells.Item(26, 5).Value = "  +1.35%  "

$ws.Cells.Item(27, 5).Value = "  +0.69%  "

$ws.Cells.Item(29, 5).Value = "  +0.01%  "

$ws.Cells.Item(30, 5).Value = "  +0.69%  "

$ws.Cells.Item(31, 5).Value = "  +1.50%  "

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.24"
$cell.ClearFormats()
$ws.Cells.Item(32, 5).Value = "  -0.04%  "

$ws.Cells.Item(33, 5).Value = "  -0.71%  "

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.398.49"
$cell.ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -4.02%  "

$ws.Cells.Item(35, 5).Value = "  -1.26%  "

$ws.Cells.Item(36, 5).Value = "  -5.54%  "

$ws.Cells.Item(37, 5).Value = "  +0.16%  "

$ws.Cells.Item(38, 5).Value = "  +0.79%  "

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.54"
$cell.ClearFormats()
$ws.Cells.Item(39, 5).Value = "  +8.65%  "

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.541"
$cell.ClearFormats()
$ws.Cells.Item(40, 5).Value = "  +0.39%  "

$ws.Cells.Item(41, 5).Value = "  +0.03%  "

$ws.Cells.Item(42, 5).Value = "  -2.57%  "

$ws.Cells.Item(43, 5).Value = "  +8.85%  "

$ws.Cells.Item(44, 5).Value = "  -0.08%  "

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.985"
$cell.ClearFormats()
$ws.Cells.Item(45, 5).Value = "  +0.28%  "

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "64.45"
$cell.ClearFormats()
$ws.Cells.Item(46, 5).Value = "  +0.08%  "

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.734.38"
$cell.ClearFormats()
$ws.Cells.Item(47, 5).Value = "  +2.06%  "

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "87.70"
$cell.ClearFormats()
$ws.Cells.Item(48, 5).Value = "  +2.15%  "

$ws.Cells.Item(49, 5).Value = "  +0.21%  "

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0₆0103"
$cell.ClearFormats()
$ws.Cells.Item(50, 5).Value = "  +0.32%  "

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0526"
$cell.ClearFormats()
$ws.Cells.Item(51, 5).Value = "  +0.34%  "
